$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (M4, N4) mirroring the style of K4/L4
# ("Multivalued" / "Unique") and set their labels to "Pattern" /
# "Pattern Type".
$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4:N4").PasteSpecial(-4122) | Out-Null

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Match the active selection recorded in the saved file.
$ws.Range("M4:N4").Select() | Out-Null
